# Update LOTOMANIA results with the latest Caixa lottery draws (concursos
# 2851-2856), appended right after the existing last row (2850, row 344).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("LOTOMANIA")

# Each inner array is: row number, Concurso, Bola1..Bola20
$newRows = @(
    ,@(345,2851,2,8,15,20,24,25,30,34,37,40,43,51,60,62,67,77,81,85,87,94)
    ,@(346,2852,1,2,4,6,17,20,21,25,33,37,38,43,53,55,67,75,83,90,94,99)
    ,@(347,2853,11,27,37,41,43,52,53,57,60,62,66,68,71,73,77,79,82,87,90,95)
    ,@(348,2854,2,4,5,7,18,25,26,28,37,38,40,42,50,54,57,68,73,76,88,89)
    ,@(349,2855,3,8,12,15,20,29,42,43,45,47,50,53,62,67,73,75,81,86,87,96)
    ,@(350,2856,0,2,4,8,9,10,12,14,15,42,43,53,67,69,72,79,85,93,98,99)
)

foreach ($rowData in $newRows) {
    $r = $rowData[0]
    for ($col = 1; $col -le 21; $col++) {
        $ws.Cells.Item($r, $col).Value = $rowData[$col]
    }
}

# Match the new selection left by the edit: the freshly entered block is
# selected with B345 as the active cell.
[void]$ws.Range("B345:U350").Select()

$win = $excel.ActiveWindow
$win.ScrollRow = 333
